# Reorder names in the "Recorded By" column (G) to reflect an updated
# recording/sync order:
#   "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#   "System, admin@admin.com"             -> "admin@admin.com, System"
#   "dnasr281@gmail.com, admin@admin.com"  -> "admin@admin.com, dnasr281@gmail.com"
#
# Entries such as "backup@backdoor.com, System" (and any 1- or 3-entry
# lists) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    $parts = $val -split ', '

    if ($parts.Count -eq 2) {
        $first = $parts[0]
        $second = $parts[1]

        $swap = $false
        if (($first -eq 'System') -and ($second -ne 'backup@backdoor.com')) {
            $swap = $true
        }
        elseif (($first -eq 'dnasr281@gmail.com') -and ($second -eq 'admin@admin.com')) {
            $swap = $true
        }

        if ($swap) {
            $cell.Value2 = $second + ', ' + $first
        }
    }
}
